# Update odds data on Sheet1 to reflect latest FlashScore values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Row 3 updates
$ws.Range("O3").Value = 1.5
$ws.Range("P3").Value = 2.5

# Row 6 updates
$ws.Range("G6").Value = 2.4
$ws.Range("H6").Value = 3.3
$ws.Range("I6").Value = 2.9
$ws.Range("O6").Value = 1.25
$ws.Range("P6").Value = 3.75
$ws.Range("Q6").Value = 1.9
$ws.Range("R6").Value = 1.95
$ws.Range("S6").Value = 1.4
$ws.Range("T6").Value = 2.75
$ws.Range("U6").Value = 1.7
$ws.Range("V6").Value = 2.05
$ws.Range("AG6").Value = 10
$ws.Range("AM6").Value = 201
$ws.Range("AT6").Value = 2.75
